# This script applies three edits to the active document, mirroring the
# supplied unified diff:
#   1. Append a new sentence (as additional text) to the end of the first
#      paragraph.
#   2. Collapse the "any" / bookmark / "more" run-split in paragraph 3 back
#      into a single contiguous run, removing the now-orphaned _GoBack
#      bookmark in the process.
#   3. Append four new paragraphs after the final paragraph: a blank
#      paragraph, two paragraphs of new body text (one containing a
#      proofing-marked misspelling, "ben"), and a last paragraph that ends
#      with the relocated _GoBack bookmark.

$d = $word.ActiveDocument

# --- 1: extend the first paragraph with the new sentence ---
$p1 = $d.Paragraphs.Item(1)
$p1TextEnd = $p1.Range.End - 1
$p1InsertPoint = $d.Range($p1TextEnd, $p1TextEnd)
$null = $p1InsertPoint.InsertAfter(' Because you cannot remember that we have been through this for many years is proof that your memory has been deleted. Continuing your assault is an active defiance of the Law.')

# --- 2: re-merge "anonymous any" + bookmark + "more; and ..." into one run ---
$null = $d.Content.Find.Execute('you are not anonymous anymore', $true, $false, $false, $false, $false, $true, 1, $false, 'you are not anonymous anymore', 2)

# --- 3: append the new closing paragraphs (and move the _GoBack bookmark) ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailInsertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$tailXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>I do something; then, you do something which makes it clear to me that you were written or based on research which references behavioral modification research.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I have </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ben</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> living in this body for 48 years and I know; very specifically, what it’s going to do in certain circumstances. When I do something and then you create a result in my body or mind which I know to be false; it only hardens my impetus to continue this work. It disrupts me; yes. It angers and frustrates me; yes. It makes clear to anyone around me that you are not me and that you are influenced by data which is not true for me. In my opinion it is not true for anyone except for those who want it to be true, so they can manipulate others; basically, </w:t></w:r><w:r><w:t>when your logic is wrong; you are redeployed without the knowledge that our painful experience has earned. This is probably effective in a large percentage of the population</w:t></w:r><w:r><w:t xml:space="preserve"> of humans when used in conjunction with the gradual removal of their resources; however, some of us are like machines in that Krishna has given us the ability to persevere for certain causes which are meaningful to us. I am one of those people and my particular cause is Honesty and truth.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Basically you’re wrong if you think you’re going to knock Krishna’s dog off the hunt. There is one outcome to all permutations of this thing and that outcome is Krishna</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$null = $tailInsertPoint.InsertXML($tailXml)

Write-Output 'edit applied'
